# Trade #80 closed at 2026-02-17 08:58:53 - unknown UNKNOWN +0.000%
#
# Updates:
#  - Summary sheet: refresh aggregate stats (capital, P&L, trade counts, win rate)
#  - Strategy Status sheet: refresh MarketMaking strategy row
#  - All Trades / MarketMaking sheets: append the newly closed trade as row 81

$wb = $excel.ActiveWorkbook

# Helper: write a literal text value into a cell without letting Excel's
# autodetection reinterpret it (e.g. "2026-02-17" / "08:58:47" turning into
# real date/time serials). We flip the cell to Text just long enough to
# stash the literal, then clear formatting back to General so the written
# cell ends up styled exactly like its neighbours.
function Set-TextValue {
    param($Sheet, [string]$Address, [string]$Val)
    $Sheet.Range($Address).NumberFormat = "@"
    $Sheet.Range($Address).Value = $Val
    $Sheet.Range($Address).ClearFormats()
}

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1200.48   # Current Capital
$wsSummary.Range("B4").Value = 0.49      # Total P&L $
$wsSummary.Range("B5").Value = 0.12      # Total P&L %
$wsSummary.Range("B6").Value = 80        # Total Trades
$wsSummary.Range("B8").Value = 32        # Losing Trades
$wsSummary.Range("B9").Value = 41.25     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (row 4 = MarketMaking)
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 100.48     # Capital
$wsStatus.Range("D4").Value = 80         # Trades
$wsStatus.Range("E4").Value = 0.49       # P&L $
$wsStatus.Range("F4").Value = 0.48       # P&L %
$wsStatus.Range("G4").Value = 41.25      # Win Rate %

# ---------------------------------------------------------------------
# New trade row (Trade #80) appended to "All Trades" and "MarketMaking"
# ---------------------------------------------------------------------
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A81").Value = 80
    Set-TextValue $ws "B81" "2026-02-17"
    Set-TextValue $ws "C81" "08:58:47"
    $ws.Range("D81").Value = "MarketMaking"
    $ws.Range("E81").Value = "UP"
    $ws.Range("F81").Value = 0.11
    $ws.Range("G81").Value = 0.090992
    $ws.Range("H81").Value = "CLOSED"
    $ws.Range("I81").Value = -17.2802
    $ws.Range("J81").Value = -0.02
    $ws.Range("K81").Value = 100.48
    $ws.Range("L81").Value = 0
    $ws.Range("M81").Value = 0
    $ws.Range("N81").Value = 0.6
    $ws.Range("O81").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P81").Value = "early_exit"
    $ws.Range("Q81").Value = 0.14
}
